# [MOSIP-14369] Fix: boolean values
#
# The "is_active" column (E) was storing TRUE as a formula (=TRUE()) that
# evaluates to the boolean number 1. It should instead be the literal text
# "TRUE", matching how the masterdata uploader expects the column to be
# populated. Rewrite E2:E11 as literal text "TRUE" instead of a boolean
# formula, then restore the cursor position recorded by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("E2:E11")

# Writing the literal word TRUE directly would be auto-coerced by Excel
# into a Boolean cell (same as the old =TRUE() formula). Enter it as a
# text formula instead, then collapse the formula down to its literal
# value with a values-only paste so the stored cell becomes plain text.
$rng.Formula = "=""TRUE"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null   # xlPasteValues

# Move the active selection to G10, as recorded in the saved workbook.
$ws.Range("G10").Select() | Out-Null
